$wb = $excel.ActiveWorkbook

# example2 sheet: add new cell E7 = "xxx" and update selection
$ws2 = $wb.Worksheets.Item("example2")
$ws2.Range("E7").Value = "xxx"

# Sheet3: remove the custom column width formatting (back to sheet default)
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A1:AMK1").EntireColumn.ClearFormats()

# make example2 the active/selected sheet (activeTab=1 in workbook.xml, tabSelected on sheet2)
$ws2.Activate()
$ws2.Range("G20").Select()
